$d = $word.ActiveDocument

# 1. "would join coalitions or coalitions would invite representatives." ->
#    "would take turns joining coalitions or coalitions would take turns inviting representatives."
$d.Content.Find.Execute(
    'would join coalitions or coalitions would invite representatives.',
    $true, $false, $false, $false, $false, $true, 1, $false,
    'would take turns joining coalitions or coalitions would take turns inviting representatives.',
    2) | Out-Null

# 2. Extend the end of the last paragraph with the new sentences (still inside the
#    same paragraph, before the existing _GoBack bookmark).
$d.Content.Find.Execute(
    'result in a more balanced budget. ',
    $true, $false, $false, $false, $false, $true, 1, $false,
    'result in a more balanced budget. One problem with both of these solutions is that it could lead to a infinite loop and never end, however the joining coalitions method resulted in a lot less of them. Therefore I will use that method. Representatives will take turns, negotiating with each coalition, they will then compare their utility by joining their coalition and choose the coalition which offers the best utility, or stay in the current coalition if its better utility then all offers.',
    2) | Out-Null

# 3. Add a trailing space run at the very end of that same paragraph (after the
#    bookmark), matching the final run added in the diff.
$lastPara = $d.Paragraphs.Last
$endRange = $lastPara.Range
$endRange.Collapse(0)
$endRange.InsertAfter(' ')

# 4. Append the three new paragraphs that follow.
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
$p1.Range.InsertAfter('However any infinite loops are unacceptable in the algorithm. To solve this I came up with two solutions. First, representatives could be impatient, this means that a budget in round n is better than a budget in round m if n < m. Another method is to make the coalition size less relevant is earlier turns and make it more important in later turns. The first method allows representatives to look forward to later rounds in order to make an more informed decision. However to work out the single negotiation would be factorial, and therefore unusable with a medium number of representatives. For this reason I will use the second method. With the new utility function as follows.  ')

$r2 = $p1.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Range.InsertAfter('The next decision. ')

$r3 = $p2.Range
$r3.Collapse(0)
$r3.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$p3.Range.InsertAfter(' ')
